# Applies the "first round of 3.3 files" update to
# Rebound Pol Emis per Unit Avoided CO2 Emis.xlsx
#
# Summary of changes:
#  - Data sheet header text updated from "Table 6-3" to "Table 6-2"
#  - Data sheet's embedded picture (old Table 6-3 excerpt image) removed
#  - Data sheet raw input values refreshed (years 2010-2014 -> 2015-2019,
#    updated CO2/CH4/N2O kt figures); dependent ratio formulas and the
#    RPEpUACE sheet's averaged rebound-factor formulas recalculate
#    automatically
#  - Active sheet/selection state updated: RPEpUACE becomes the active
#    (selected) tab instead of About; Data sheet's selection moves to E25

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")
$wsR     = $wb.Worksheets.Item("RPEpUACE")

# --- Data sheet: update the source excerpt label -----------------------
$wsData.Range("A1").Value = "Excerpt from Table 6-2:"

# --- Data sheet: remove the embedded picture ----------------------------
while ($wsData.Shapes.Count -gt 0) {
    $wsData.Shapes.Item(1).Delete()
}

# --- Data sheet: refresh the year headers -------------------------------
$wsData.Range("B3").Value = 2015
$wsData.Range("C3").Value = 2016
$wsData.Range("D3").Value = 2017
$wsData.Range("E3").Value = 2018
$wsData.Range("F3").Value = 2019

# --- Data sheet: refresh CO2 (kt) row -----------------------------------
$wsData.Range("B4").Value = -791695
$wsData.Range("C4").Value = -855998
$wsData.Range("D4").Value = -792046
$wsData.Range("E4").Value = -824885
$wsData.Range("F4").Value = -812695

# --- Data sheet: refresh CH4 (kt) row -----------------------------------
$wsData.Range("B5").Value = 663
$wsData.Range("C5").Value = 308
$wsData.Range("D5").Value = 614
$wsData.Range("E5").Value = 552
$wsData.Range("F5").Value = 552

# --- Data sheet: refresh N2O (kt) row -----------------------------------
$wsData.Range("B6").Value = 38
$wsData.Range("C6").Value = 18
$wsData.Range("D6").Value = 36
$wsData.Range("E6").Value = 32
$wsData.Range("F6").Value = 32

# Rows 8:9 (CH4/CO2 and N2O/CO2 ratios) and RPEpUACE!B11:B12 (the averaged
# rebound emission factors) are formulas referencing the cells above, so
# they recompute automatically once the inputs change.

# --- Selection / active-tab bookkeeping ---------------------------------
# Data sheet's remembered selection moves from F7 to E25.
$wsData.Range("E25").Select()

# RPEpUACE becomes the active/selected sheet (was About); its own
# remembered selection (B12) is unchanged.
$wsR.Activate()
$wsR.Range("B12").Select()
